$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1) Change the LUCAS row (account 005242683) in place to CARLA's data.
$lucasRow = $ws.Cells.Find("005242683").Row
$ws.Cells.Item($lucasRow, 1).NumberFormat = "@"
$ws.Cells.Item($lucasRow, 1).Value = "004643153"
$ws.Cells.Item($lucasRow, 2).Value = "CARLA"
$ws.Cells.Item($lucasRow, 3).Value = 5372.43

# 2) Delete the rows that disappear entirely (bottom-most first so row
#    numbers found earlier stay valid).
$danielaDupRow = $ws.Cells.Find("004001621").Row
$ws.Rows.Item($danielaDupRow).Delete()

$priscillaRow = $ws.Cells.Find("004224284").Row
$ws.Rows.Item($priscillaRow).Delete()

$lagoRow = $ws.Cells.Find("001882235").Row
$ws.Rows.Item($lagoRow + 2).Delete()
$ws.Rows.Item($lagoRow + 1).Delete()
$ws.Rows.Item($lagoRow).Delete()

# 3) Insert the new MARIAH row right before CAIO (004512434).
$caioRow = $ws.Cells.Find("004512434").Row
$ws.Rows.Item($caioRow).Insert()
$ws.Cells.Item($caioRow, 1).NumberFormat = "@"
$ws.Cells.Item($caioRow, 1).Value = "004242237"
$ws.Cells.Item($caioRow, 2).Value = "MARIAH"
$ws.Cells.Item($caioRow, 3).Value = 2034.07

# 4) Insert the new DANIELA row right before FRANCISCO (004567324), at the
#    top of the data (just under the header).
$franciscoRow = $ws.Cells.Find("004567324").Row
$ws.Rows.Item($franciscoRow).Insert()
$ws.Cells.Item($franciscoRow, 1).NumberFormat = "@"
$ws.Cells.Item($franciscoRow, 1).Value = "004001621"
$ws.Cells.Item($franciscoRow, 2).Value = "DANIELA"
$ws.Cells.Item($franciscoRow, 3).Value = 100037.58
